$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "294.28"
Set-TextValue "E2" "-2.58%"
Set-TextValue "D3" "31.13"
Set-TextValue "E3" "-3.13%"
Set-TextValue "D4" "4.965"
Set-TextValue "E4" "-1.21%"
Set-TextValue "D5" "0.07324"
Set-TextValue "E5" "-7.10%"
Set-TextValue "D6" "1.789"
Set-TextValue "E6" "-14.76%"
Set-TextValue "D7" "7.661"
Set-TextValue "E7" "-1.92%"
Set-TextValue "D8" "3.762"
Set-TextValue "E8" "-0.76%"
Set-TextValue "D9" "0.9071"
Set-TextValue "E9" "-2.07%"
Set-TextValue "D10" "0.1649"
Set-TextValue "E10" "-5.91%"
Set-TextValue "D11" "0.07561"
Set-TextValue "E11" "-5.08%"
Set-TextValue "D12" "0.08117"
Set-TextValue "E12" "-7.84%"
Set-TextValue "D13" "0.02990"
Set-TextValue "E13" "-4.12%"
Set-TextValue "D14" "0.09989"
Set-TextValue "E14" "-0.44%"
Set-TextValue "D15" "0.001498"
Set-TextValue "E15" "-0.60%"
Set-TextValue "D16" "0.005611"
Set-TextValue "E16" "-4.86%"
Set-TextValue "D17" "3.459"
Set-TextValue "E17" "-0.14%"
Set-TextValue "D18" "2.095"
Set-TextValue "E18" "-8.08%"
Set-TextValue "D19" "0.3271"
Set-TextValue "E19" "-0.67%"
Set-TextValue "D20" "0.1307"
Set-TextValue "E20" "1.45%"
Set-TextValue "D21" "4.364"
Set-TextValue "E21" "4.81%"
Set-TextValue "D22" "0.2001"
Set-TextValue "E22" "11.90%"
Set-TextValue "D23" "0.04472"
Set-TextValue "E23" "-2.62%"
Set-TextValue "D24" "0.001225"
Set-TextValue "E24" "-0.80%"
Set-TextValue "D25" "0.004040"
Set-TextValue "E25" "-10.46%"
Set-TextValue "D26" "0.0001251"
Set-TextValue "E26" "0.31%"
Set-TextValue "D39" "0.01644"
Set-TextValue "E39" "-5.56%"
Set-TextValue "D40" "0.04393"
Set-TextValue "E40" "-7.79%"
Set-TextValue "D41" "0.007363"
Set-TextValue "E41" "-0.10%"
Set-TextValue "D42" "0.1318"
Set-TextValue "E42" "-3.68%"
Set-TextValue "D43" "0.002063"
Set-TextValue "E43" "-3.37%"
Set-TextValue "D44" "0.01121"
Set-TextValue "E44" "4.11%"
Set-TextValue "D45" "0.00005992"
Set-TextValue "E45" "-1.27%"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "E46" "0.31%"
Set-TextValue "D47" "2.129"
Set-TextValue "E47" "159.48%"
Set-TextValue "D48" "0.002402"
Set-TextValue "E48" "-29.21%"
Set-TextValue "D49" "0.00002102"
Set-TextValue "E49" "0.31%"
Set-TextValue "D50" "0.0002002"
Set-TextValue "E50" "0.31%"
